$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")
$ws.Range("A2").Value = "test"
